$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values in columns I/J on rows 2-4 (Asset Store things used / Libraries used)
$ws.Range("I2").Value = "Torch"
$ws.Range("J2").Value = "FPS"
$ws.Range("I3").Value = "Rain AI"
$ws.Range("I4").Value = "Probuilder"

# New log rows 36-37. Columns B/C hold dates that are stored as plain text
# (not date serials) in the source file, so force text formatting before
# entering the values and then clear the formatting back to the default
# style, matching how the workbook already stores other text-like dates
# (e.g. B18/B19) with no explicit style index.
$datesRange = $ws.Range("B36:C37")
$datesRange.NumberFormat = "@"

$ws.Range("A36").Value = "More things added"
$ws.Range("B36").Value = "04/13/2017"
$ws.Range("C36").Value = "04/20/2017"
$ws.Range("E36").Value = "Added a new stage, fixed AI, where adding more AI will require their own scripts"

$ws.Range("A37").Value = "Objectives"
$ws.Range("B37").Value = "04/13/2017"
$ws.Range("C37").Value = "04/20/2017"
$ws.Range("E37").Value = "Ability to swap to another scene if player makes it to exit."

$datesRange.ClearFormats()

$ws.Range("J2").Select()
